$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F
$ws.Range("F1").Value = "sd_temp"

# Add sd_temp values for rows 2-17
$values = @(
    4.10186018688312,
    3.41625508938201,
    2.65425669226604,
    2.77630646996135,
    3.33303520919959,
    4.44635981800796,
    2.21413712731855,
    2.71130383327513,
    2.52225280536793,
    5.35271320395833,
    2.49821592226367,
    4.15503998807876,
    4.479966145188,
    3.07316679817093,
    3.071419894473,
    4.08808678289712
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 6).Value = $val
    $row = $row + 1
}
